$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.387899999999991
$ws.Range("D4").Value = -7.2623
$ws.Range("C7").Value = -12.8329
$ws.Range("A8").Value = -22.20610000000001
$ws.Range("A10").Value = -21.59719999999998
$ws.Range("D11").Value = -7.8477
$ws.Range("A12").Value = -21.5156
$ws.Range("C14").Value = -13.3073
$ws.Range("D14").Value = -8.042799999999991
$ws.Range("C15").Value = -14.10369999999998
$ws.Range("A18").Value = -22.25660000000003
$ws.Range("C18").Value = -13.256
$ws.Range("D18").Value = -8.112199999999991
$ws.Range("D19").Value = -8.976099999999994
$ws.Range("C20").Value = -11.823
$ws.Range("D21").Value = -8.731500000000002
$ws.Range("A25").Value = -21.34689999999998
$ws.Range("D27").Value = -8.537000000000006
$ws.Range("C29").Value = -11.9426
$ws.Range("C30").Value = -12.45909999999999
$ws.Range("C31").Value = -12.7661
$ws.Range("D31").Value = -8.764800000000005
$ws.Range("C35").Value = -11.629
$ws.Range("A37").Value = -20.59710000000002
$ws.Range("D38").Value = -8.413200000000005
$ws.Range("C40").Value = -13.5187
$ws.Range("D42").Value = -8.753599999999993
$ws.Range("C44").Value = -13.48179999999999
$ws.Range("D44").Value = -7.8164
$ws.Range("D47").Value = -7.750199999999999
$ws.Range("C50").Value = -13.3806
$ws.Range("C54").Value = -13.261
$ws.Range("A55").Value = -22.4234
$ws.Range("D56").Value = -7.980399999999999
$ws.Range("D58").Value = -8.45249999999999
$ws.Range("D65").Value = -7.952899999999996
$ws.Range("A68").Value = -21.65899999999999
$ws.Range("C68").Value = -11.7629
$ws.Range("D73").Value = -7.936999999999997
$ws.Range("C76").Value = -12.3662
$ws.Range("A77").Value = -20.77909999999999
$ws.Range("A78").Value = -20.85209999999999
$ws.Range("A79").Value = -20.87849999999998
$ws.Range("A80").Value = -20.65150000000001
$ws.Range("A81").Value = -21.86990000000001
$ws.Range("A82").Value = -21.8503
$ws.Range("A84").Value = -21.9498
$ws.Range("C87").Value = -13.34889999999998
$ws.Range("C88").Value = -12.6981
$ws.Range("D90").Value = -7.972900000000005
$ws.Range("C92").Value = -11.4213
$ws.Range("D92").Value = -6.671300000000002
$ws.Range("D94").Value = -6.867999999999999
$ws.Range("D95").Value = -7.925599999999999
$ws.Range("C96").Value = -12.55500000000001
$ws.Range("C98").Value = -11.7176
$ws.Range("A101").Value = -21.27989999999999
$ws.Range("C101").Value = -12.95620000000001
$ws.Range("D101").Value = -8.197200000000002
$ws.Range("A102").Value = -19.93049999999999
$ws.Range("C102").Value = -13.26260000000001
